$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Inscricoes")

$ws.Range("F2").Value = 39
$ws.Range("H2").Value = 51
$ws.Range("F3").Value = 16
$ws.Range("H3").Value = 20
$ws.Range("E4").Value = 27
$ws.Range("E12").Value = 5
$ws.Range("F12").Value = 2
$ws.Range("H12").Value = 2
$ws.Range("E15").Value = 169
$ws.Range("F15").Value = 95
$ws.Range("H15").Value = 136
$ws.Range("E17").Value = 134
$ws.Range("F17").Value = 71
$ws.Range("H17").Value = 103
$ws.Range("E18").Value = 124
$ws.Range("F18").Value = 58
$ws.Range("H18").Value = 94
$ws.Range("E19").Value = 68
$ws.Range("F19").Value = 42
$ws.Range("H19").Value = 55
$ws.Range("F23").Value = 3
$ws.Range("H23").Value = 5
$ws.Range("E24").Value = 27
$ws.Range("F25").Value = 13
$ws.Range("H25").Value = 21
$ws.Range("F26").Value = 17
$ws.Range("H26").Value = 27
$ws.Range("E27").Value = 12
$ws.Range("F29").Value = 12
$ws.Range("H29").Value = 15
$ws.Range("E31").Value = 2
$ws.Range("F31").Value = 1
$ws.Range("H31").Value = 2
$ws.Range("F34").Value = 10
$ws.Range("H34").Value = 13
$ws.Range("F35").Value = 6
$ws.Range("H35").Value = 7
$ws.Range("F36").Value = 58
$ws.Range("H36").Value = 90
$ws.Range("F37").Value = 37
$ws.Range("H37").Value = 49
$ws.Range("E40").Value = 26
$ws.Range("F40").Value = 17
$ws.Range("H40").Value = 20
$ws.Range("F41").Value = 23
$ws.Range("H41").Value = 34
$ws.Range("F42").Value = 24
$ws.Range("H42").Value = 33
$ws.Range("F44").Value = 17
$ws.Range("H44").Value = 27
$ws.Range("E47").Value = 64
$ws.Range("F47").Value = 40
$ws.Range("H47").Value = 51
$ws.Range("E48").Value = 36
$ws.Range("F48").Value = 24
$ws.Range("H48").Value = 29
$ws.Range("E50").Value = 29
$ws.Range("F58").Value = 3
$ws.Range("H58").Value = 3
$ws.Range("F61").Value = 14
$ws.Range("H61").Value = 24
$ws.Range("F62").Value = 16
$ws.Range("H62").Value = 30
$ws.Range("E65").Value = 38
$ws.Range("F65").Value = 13
$ws.Range("H65").Value = 26
$ws.Range("F66").Value = 25
$ws.Range("H66").Value = 33
$ws.Range("E68").Value = 20
$ws.Range("F71").Value = 22
$ws.Range("H71").Value = 32
$ws.Range("E72").Value = 47
$ws.Range("F72").Value = 28
$ws.Range("H72").Value = 39
$ws.Range("F73").Value = 13
$ws.Range("H73").Value = 25
$ws.Range("E75").Value = 17
$ws.Range("E79").Value = 43
$ws.Range("F79").Value = 22
$ws.Range("H79").Value = 34
$ws.Range("F80").Value = 17
$ws.Range("H80").Value = 29
$ws.Range("E84").Value = 6
$ws.Range("E88").Value = 29
$ws.Range("F89").Value = 22
$ws.Range("H89").Value = 29
